$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Qminus1)
$ws.Range("B2").Value = -0.01879798516539222
$ws.Range("C2").Value = 0.3559641005099519
$ws.Range("D2").Value = 0.2241929700663822
$ws.Range("E2").Value = 0.4734902006022746
$ws.Range("F2").Value = 0.4812049593648728

# Row 3 (Q0)
$ws.Range("B3").Value = 0.1997949562989836
$ws.Range("C3").Value = 0.5133874392789209
$ws.Range("D3").Value = 0.4202257940472184
$ws.Range("E3").Value = 0.6482482503232989
$ws.Range("F3").Value = 0.6199282290422755
$ws.Range("G3").Value = 96

# Row 4 (Q1)
$ws.Range("B4").Value = 0.1942849321833074
$ws.Range("C4").Value = 0.5188225552120047
$ws.Range("D4").Value = 0.3557134039314317
$ws.Range("E4").Value = 0.5964171391999326
$ws.Range("F4").Value = 0.5701163892217513
$ws.Range("G4").Value = 46
